$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = 8.720000267028809
$ws.Range("E2").Value = 6.5
$ws.Range("F2").Value = 10.80000019073486
$ws.Range("G2").Value = 5.75
$ws.Range("H2").Value = 375150890
$ws.Range("I2").Value = "SOUN"
$ws.Range("D3").Value = 8.720000267028809
$ws.Range("E3").Value = 6.5
$ws.Range("F3").Value = 10.80000019073486
$ws.Range("G3").Value = 5.75
$ws.Range("H3").Value = 375150890
$ws.Range("I3").Value = "SOUN"
$ws.Range("D4").Value = 8.720000267028809
$ws.Range("E4").Value = 6.5
$ws.Range("F4").Value = 10.80000019073486
$ws.Range("G4").Value = 5.75
$ws.Range("H4").Value = 375150890
$ws.Range("I4").Value = "SOUN"
$ws.Range("D5").Value = 8.720000267028809
$ws.Range("E5").Value = 6.5
$ws.Range("F5").Value = 10.80000019073486
$ws.Range("G5").Value = 5.75
$ws.Range("H5").Value = 375150890
$ws.Range("I5").Value = "SOUN"
$ws.Range("D6").Value = 8.720000267028809
$ws.Range("E6").Value = 6.5
$ws.Range("F6").Value = 10.80000019073486
$ws.Range("G6").Value = 5.75
$ws.Range("H6").Value = 375150890
$ws.Range("I6").Value = "SOUN"
$ws.Range("D7").Value = 8.720000267028809
$ws.Range("E7").Value = 6.5
$ws.Range("F7").Value = 10.80000019073486
$ws.Range("G7").Value = 5.75
$ws.Range("H7").Value = 375150890
$ws.Range("I7").Value = "SOUN"
$ws.Range("D8").Value = 2.579999923706055
$ws.Range("E8").Value = 3.5
$ws.Range("F8").Value = 4.53000020980835
$ws.Range("G8").Value = 2.400000095367432
$ws.Range("H8").Value = 375150890
$ws.Range("I8").Value = "SOUN"
$ws.Range("D9").Value = 3.369999885559082
$ws.Range("E9").Value = 2.710000038146973
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 2.5
$ws.Range("H9").Value = 375150890
$ws.Range("I9").Value = "SOUN"
$ws.Range("D10").Value = 1.75
$ws.Range("E10").Value = 1.990000009536743
$ws.Range("F10").Value = 2.490000009536743
$ws.Range("G10").Value = 0.9700000286102296
$ws.Range("H10").Value = 375150890
$ws.Range("I10").Value = "SOUN"
$ws.Range("D11").Value = 2.759999990463257
$ws.Range("E11").Value = 2.660000085830688
$ws.Range("F11").Value = 3.345000028610229
$ws.Range("G11").Value = 2.224999904632568
$ws.Range("H11").Value = 375150890
$ws.Range("I11").Value = "SOUN"
$ws.Range("D12").Value = 4.590000152587891
$ws.Range("E12").Value = 2.329999923706055
$ws.Range("F12").Value = 4.748000144958496
$ws.Range("G12").Value = 2.075000047683716
$ws.Range("H12").Value = 375150890
$ws.Range("I12").Value = "SOUN"
$ws.Range("D13").Value = 2.069999933242798
$ws.Range("E13").Value = 1.590000033378601
$ws.Range("F13").Value = 2.089999914169312
$ws.Range("G13").Value = 1.490000009536743
$ws.Range("H13").Value = 375150890
$ws.Range("I13").Value = "SOUN"
$ws.Range("D14").Value = 2.130000114440918
$ws.Range("E14").Value = 1.659999966621399
$ws.Range("F14").Value = 2.176000118255615
$ws.Range("G14").Value = 1.649999976158142
$ws.Range("H14").Value = 375150890
$ws.Range("I14").Value = "SOUN"
$ws.Range("D15").Value = 5.909999847412109
$ws.Range("E15").Value = 4.239999771118164
$ws.Range("F15").Value = 5.940000057220459
$ws.Range("G15").Value = 3.5
$ws.Range("H15").Value = 375150890
$ws.Range("I15").Value = "SOUN"
$ws.Range("D16").Value = 3.980000019073486
$ws.Range("E16").Value = 5.090000152587891
$ws.Range("F16").Value = 6.449999809265137
$ws.Range("G16").Value = 3.819999933242798
$ws.Range("H16").Value = 375150890
$ws.Range("I16").Value = "SOUN"
$ws.Range("D17").Value = 4.690000057220459
$ws.Range("E17").Value = 5.03000020980835
$ws.Range("F17").Value = 6.25
$ws.Range("G17").Value = 4.449999809265137
$ws.Range("H17").Value = 375150890
$ws.Range("I17").Value = "SOUN"
$ws.Range("D18").Value = 20.47500038146973
$ws.Range("E18").Value = 14.14999961853027
$ws.Range("F18").Value = 22.85000038146973
$ws.Range("G18").Value = 12.40999984741211
$ws.Range("H18").Value = 375150890
$ws.Range("I18").Value = "SOUN"
$ws.Range("D19").Value = 8.270000457763672
$ws.Range("E19").Value = 9.289999961853027
$ws.Range("F19").Value = 9.970000267028809
$ws.Range("G19").Value = 6.519999980926514
$ws.Range("H19").Value = 375150890
$ws.Range("I19").Value = "SOUN"
$ws.Range("D20").Value = 10.52499961853027
$ws.Range("E20").Value = 10.32999992370606
$ws.Range("F20").Value = 13.55500030517578
$ws.Range("G20").Value = 9.890000343322754
$ws.Range("H20").Value = 375150890
$ws.Range("I20").Value = "SOUN"
